$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => [D value (or $null to leave unchanged), E value]
$updates = @{
    2  = @("30.540.78", "  +1.44%  ")
    3  = @("1.881.46", "  +1.35%  ")
    4  = @($null, "  -0.04%  ")
    5  = @("248.03", "  +6.53%  ")
    7  = @("0.4754", "  +1.69%  ")
    8  = @("0.2924", "  +3.43%  ")
    9  = @("0.06523", "  +1.25%  ")
    10 = @("21.98", "  +6.11%  ")
    11 = @("98.04", "  +5.07%  ")
    12 = @("0.07720", "  +0.70%  ")
    13 = @("0.7392", "  +9.44%  ")
    14 = @("1.884.85", "  +1.44%  ")
    15 = @("5.144", "  +2.32%  ")
    16 = @("274.01", "  +2.83%  ")
    17 = @("30.528.26", "  +1.45%  ")
    18 = @("13.47", "  +1.49%  ")
    19 = @("0.000007566", "  +1.00%  ")
    20 = @($null, "  -0.03%  ")
    21 = @("2.133.07", "  +1.94%  ")
    22 = @($null, "  -0.02%  ")
    23 = @("5.249", $null)
    24 = @("6.197", "  +2.07%  ")
    25 = @("9.284", "  +0.54%  ")
    26 = @("163.53", "  -1.02%  ")
    27 = @("18.89", "  +1.83%  ")
    28 = @("1.943", "  +3.88%  ")
    29 = @($null, "  +2.96%  ")
    30 = @($null, "  -0.55%  ")
    31 = @("1.520", "  +5.19%  ")
    32 = @("4.327", "  +3.62%  ")
    33 = @("4.109", "  +3.96%  ")
    34 = @("0.04822", "  +4.35%  ")
    35 = @($null, "  +2.32%  ")
    36 = @("0.7015", "  +3.25%  ")
    37 = @("2.714", "  +0.04%  ")
    38 = @("0.01869", "  +3.28%  ")
    39 = @($null, "  +1.69%  ")
    40 = @("6.321", "  +0.73%  ")
    41 = @("1.994", "  +7.03%  ")
    42 = @("71.33", "  +1.85%  ")
    43 = @("0.4218", "  +4.92%  ")
    44 = @("0.8418", "  +1.90%  ")
    46 = @("102.71", "  +0.82%  ")
    47 = @("9.347", "  +1.77%  ")
    48 = @($null, "  +3.31%  ")
    49 = @("35.61", "  +4.91%  ")
    50 = @("916.60", "  -0.26%  ")
    51 = @("0.3891", "  +4.39%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]
    if ($null -ne $dVal) {
        $cell = $ws.Cells.Item($row, 4)
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $dVal
        $cell.Style = $origStyle
    }
    if ($null -ne $eVal) {
        $cell = $ws.Cells.Item($row, 5)
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $eVal
        $cell.Style = $origStyle
    }
}
